# feat(CWL): mod integrity check
# Adds three new rows (60-62) describing the new "missing mods" warning
# dialog strings to the EN language sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Formatting first (copy/paste of formats does not touch shared
# strings, so this can happen in whatever order is convenient).
# ---------------------------------------------------------------------

# Column A ("id") on all three rows - reuse the existing key-column
# formatting (font 5: 15.8pt Cascadia Code, FFC000) from A59.
$ws.Cells.Item(59, 1).Copy()
$ws.Cells.Item(60, 1).PasteSpecial(-4122)
$ws.Cells.Item(59, 1).Copy()
$ws.Cells.Item(61, 1).PasteSpecial(-4122)
$ws.Cells.Item(59, 1).Copy()
$ws.Cells.Item(62, 1).PasteSpecial(-4122)

# Column D ("text") on all three rows - reuse the existing EN body
# formatting (font 2, wrap) from D59.
$ws.Cells.Item(59, 4).Copy()
$ws.Cells.Item(60, 4).PasteSpecial(-4122)
$ws.Cells.Item(59, 4).Copy()
$ws.Cells.Item(61, 4).PasteSpecial(-4122)
$ws.Cells.Item(59, 4).Copy()
$ws.Cells.Item(62, 4).PasteSpecial(-4122)

# Column C ("text_JP"):
#  - C60/C61 reuse the long-body wrap formatting (font 6, from C46),
#    then get renamed onto the new 15.8pt SimSun (宋体) font.
#  - C62 uses the plain/default-ish 11pt font (no theme color).
$ws.Cells.Item(46, 3).Copy()
$ws.Cells.Item(60, 3).PasteSpecial(-4122)
$ws.Cells.Item(60, 3).Font.Name = "宋体"

$ws.Cells.Item(46, 3).Copy()
$ws.Cells.Item(61, 3).PasteSpecial(-4122)
$ws.Cells.Item(61, 3).Font.Name = "宋体"

$excel.CutCopyMode = 0

# Row heights.
$ws.Rows.Item(60).RowHeight = 46.5
$ws.Rows.Item(61).RowHeight = 23.25
$ws.Rows.Item(62).RowHeight = 23.25

# ---------------------------------------------------------------------
# Now populate values. The relative order of these assignments decides
# the order new entries are appended to the shared-string table, so it
# is kept deliberate: column A top-to-bottom, then column D (61,62,60),
# then column C (61,62,60) - matching how the sheet was authored.
# ---------------------------------------------------------------------

$ws.Cells.Item(60, 1).Value = "cwl_warn_missing_mods"
$ws.Cells.Item(61, 1).Value = "cwl_warn_missing_mods_yes"
$ws.Cells.Item(62, 1).Value = "cwl_warn_missing_mods_no"

$ws.Cells.Item(61, 4).Value = "Quit Without Saving"
$ws.Cells.Item(62, 4).Value = "Continue Playing"
$ws.Cells.Item(60, 4).Value = "Mods missing from current save:`n{0}"

$ws.Cells.Item(61, 3).Value = "セーブせずに終了"

$ws.Cells.Item(62, 3).Value = "プレイを続ける"
$ws.Cells.Item(62, 3).Font.ColorIndex = -4105

$run1 = "現在のセーブから欠落している"
$run2 = "MOD"
$run3 = "：`n{0}"
$c60c = $ws.Cells.Item(60, 3)
$c60c.Value = $run1 + $run2 + $run3
$len1 = $run1.Length
$len2 = $run2.Length
$len3 = $run3.Length
$r2 = $c60c.Characters($len1 + 1, $len2)
$r2.Font.Name = "Cascadia Code"
$r2.Font.ColorIndex = -4105
$r3 = $c60c.Characters($len1 + $len2 + 1, $len3)
$r3.Font.Name = "宋体"
$r3.Font.ColorIndex = -4105

# ---------------------------------------------------------------------
# Restore view state to approximate author's saved selection/scroll.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$win.ScrollColumn = 1
$ws.Range("C65").Select()

$excel.CutCopyMode = 0
